$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Country" (sheet1): update the single data row.
# ---------------------------------------------------------------------------
$wsCountry = $wb.Worksheets.Item("Country")
$wsCountry.Range("A2").Value = 18.64
$wsCountry.Range("B2").Value = 0.7568

# ---------------------------------------------------------------------------
# Sheet "States" (sheet2): the state rows were re-sorted by CONCERN.HIGH
# (descending) against the refreshed metrics, and the three all-zero rows
# that no longer make the cut (old Gujarat/Tamil Nadu/Karnataka rows) were
# dropped, shrinking the table from 30 to 27 states.
# Easiest/most reliable way to land on the exact target layout: wipe the
# existing data block and rewrite it row-by-row in the new order.
# ---------------------------------------------------------------------------
$wsStates = $wb.Worksheets.Item("States")
$wsStates.Range("A2:E31").ClearContents()

$statesData = @"
IN-BR|Bihar|78.95|Bihar|0
IN-MN|Manipur|50|Manipur|0
IN-UP|Uttar Pradesh|49.33|Uttar Pradesh|0
IN-JH|Jharkhand|45.83|Jharkhand|0
IN-MZ|Mizoram|36.36|Mizoram|-20
IN-NL|Nagaland|36.36|Nagaland|33.3333
IN-AR|Arunachal Pradesh|32|Arunachal Pradesh|33.3333
IN-PB|Punjab|31.82|Punjab|0
IN-PY|Puducherry|25|Puducherry|0
IN-JK|Jammu and Kashmir|22.73|Jammu and Kashmir|0
IN-ML|Meghalaya|18.18|Meghalaya|0
IN-TR|Tripura|12.5|Tripura|0
IN-TS|Telangana|12.12|Telangana|0
IN-AS|Assam|12.12|Assam|0
IN-HR|Haryana|9.09|Haryana|0
IN-DL|Delhi|9.09|Delhi|0
IN-HP|Himachal Pradesh|8.33|Himachal Pradesh|0
IN-MP|Madhya Pradesh|7.69|Madhya Pradesh|0
IN-OR|Odisha|6.67|Odisha|0
IN-CT|Chhattisgarh|3.7|Chhattisgarh|0
IN-DD|Daman and Diu|0|Daman and Diu|-100
IN-RJ|Rajasthan|0|Rajasthan|NA
IN-WB|West Bengal|0|West Bengal|NA
IN-MH|Maharashtra|0|Maharashtra|NA
IN-GJ|Gujarat|0|Gujarat|NA
IN-KA|Karnataka|0|Karnataka|NA
IN-TN|Tamil Nadu|0|Tamil Nadu|NA
"@

$rowNum = 2
foreach ($line in ($statesData -split "`n")) {
    $t = $line.Trim()
    if ($t.Length -eq 0) { continue }
    $f = $t -split '\|'

    $wsStates.Cells.Item($rowNum, 1).Value = $f[0]
    $wsStates.Cells.Item($rowNum, 2).Value = $f[1]
    $wsStates.Cells.Item($rowNum, 3).Value = [double]$f[2]
    $wsStates.Cells.Item($rowNum, 4).Value = $f[3]
    if ($f[4] -ne "NA") {
        $wsStates.Cells.Item($rowNum, 5).Value = [double]$f[4]
    }

    $rowNum++
}

# ---------------------------------------------------------------------------
# Sheet "Dark clusters" (sheet3): refreshed metrics for the three clusters;
# Ghaggar now also has a MoM figure where it previously had none.
# ---------------------------------------------------------------------------
$wsClusters = $wb.Worksheets.Item("Dark clusters")
$wsClusters.Range("B2").Value = 80
$wsClusters.Range("C2").Value = 0
$wsClusters.Range("B3").Value = 35
$wsClusters.Range("C3").Value = 1.0976
$wsClusters.Range("B4").Value = 26.96
$wsClusters.Range("C4").Value = 6.8993
